$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing weekly columns right by inserting 3 new columns before column B.
# (old B:E -> new E:H; new B,C,D are freshly inserted blank columns)
$ws.Range("B1:D1").EntireColumn.Insert()

# New latest-week header: B1 becomes the new "Jun_27" column, C1/D1 are
# additional "Jun_26" columns inserted alongside it.
$ws.Range("B1").Value = "Jun_27"
$ws.Range("C1").Value = "Jun_26"
$ws.Range("D1").Value = "Jun_26"

# Fill the newly inserted B/C/D columns for every data row with the same
# "UN" placeholder used throughout the rest of the grid.
$ws.Range("B2:D27").Value = "UN"

# Append the two new rows for the newly tracked research firms.
$ws.Range("A28").Value = "Benchmark"
$ws.Range("B28:D28").Value = "UN"

$ws.Range("A29").Value = "Evercore ISI"
$ws.Range("B29:D29").Value = "UN"
